$d = $word.ActiveDocument

# Locate the paragraph that immediately follows the one we need to
# duplicate/insert before ("Если выполнение цепочки команд завершилось...").
$rng = $d.Content
$found = $rng.Find.Execute(
    "Если выполнение цепочки команд завершилось с ошибкой",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$targetPara = $rng.Paragraphs(1)

# Insert a fresh paragraph immediately before it. It inherits the
# surrounding list/paragraph formatting (pStyle "Normal", bullet list
# numId 2 / ilvl 0, jc "both") but carries no run-level coloring.
$targetPara.Range.InsertParagraphBefore()

# After the insertion, $targetPara now resolves to the newly created
# (still empty) paragraph, so we just fill in its text.
$targetPara.Range.Text = "После разбора запроса от клиента и создания цепочки команд, эта цепочка команд передается контексту выполнения; контекст выполнения — это отдельный процесс, выполняющий каждую команду из цепочки в своем процессе. "
